$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value from EB000016 to EB000013
$ws.Range("B2").Value = "EB000013"

# Add new row 3 with the old B2 value EB000016
$ws.Range("B3").Value = "EB000016"

# Update selection to C10
$ws.Range("C10").Select()
